# Fix date errors in original data (rows mistakenly entered as 2020 dates
# instead of 2019) and restore the sheet's last-known view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the mis-typed year-2020 dates back to the intended 2019 dates ---
$ws.Range("A3").Value  = 43800   # was 44166 (2020-12-01) -> 2019-12-01
$ws.Range("A23").Value = 43812   # was 44178 (2020-12-13) -> 2019-12-13
$ws.Range("A52").Value = 43827   # was 44193 (2020-12-28) -> 2019-12-28
$ws.Range("A53").Value = 43825   # was 44191 (2020-12-26) -> 2019-12-26

# --- Restore the sheet's scroll position / active selection ---
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A56").Select()
